$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.128.67'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.13%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.847.57'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.31%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '235.76'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.94%  '
$ws.Range("E6").Value = '  -0.12%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4752'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.75%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2816'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.75%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06470'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.00%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.853.49'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.52%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07296'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.62%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '16.30'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.27%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.106'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.46%  '
$ws.Range("E14").Value = '  -1.20%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6444'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.87%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '30.076.65'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.20%  '
$ws.Range("E17").Value = '  -0.07%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.21'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.56%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007604'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.61%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.110.92'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.71%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.002'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.05%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.254'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.41%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '215.94'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +14.68%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.091'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.89%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.178'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.97%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '162.91'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.19%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.32'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.14%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.912'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.96%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.430'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.44%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09170'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.18%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.225'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.80%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.957'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.36%  '
$ws.Range("E33").Value = '  -3.82%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7386'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.26%  '
$ws.Range("E35").Value = '  +3.13%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.686'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.24%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01817'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.09%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.600'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.67%  '
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.059'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.38%  '
$ws.Range("B40").Value = 'TrustWalletToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.8995'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.84%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.908'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.29%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '106.04'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.55%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9993'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.68%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4235'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.51%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.408'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.35%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1303'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -5.82%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.548'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +9.91%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '63.82'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -6.27%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.769'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.01%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '34.14'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.16%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05678'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.58%  '
